{"js": "// Replace the three-digit x one-digit multiplication problems/answers\n// in the worksheet table with a newly generated set of problems.\n// Each old equation string is unique in the document, so a targeted\n// search-and-replace per pair is safe and precise.\nconst replacements = [\n  [\"933\u00d74=3732\", \"425\u00d77=2975\"],\n  [\"738\u00d79=6642\", \"468\u00d79=4212\"],\n  [\"161\u00d74=644\", \"766\u00d77=5362\"],\n  [\"518\u00d76=3108\", \"173\u00d78=1384\"],\n  [\"567\u00d79=5103\", \"101\u00d76=606\"],\n  [\"841\u00d79=7569\", \"173\u00d74=692\"],\n  [\"283\u00d77=1981\", \"562\u00d73=1686\"],\n  [\"359\u00d75=1795\", \"421\u00d75=2105\"],\n  [\"256\u00d74=1024\", \"605\u00d75=3025\"],\n  [\"171\u00d79=1539\", \"291\u00d75=1455\"],\n  [\"997\u00d75=4985\", \"252\u00d79=2268\"],\n  [\"600\u00d72=1200\", \"947\u00d76=5682\"],\n  [\"906\u00d78=7248\", \"897\u00d75=4485\"],\n  [\"281\u00d76=1686\", \"565\u00d79=5085\"],\n  [\"685\u00d79=6165\", \"694\u00d79=6246\"],\n  [\"378\u00d72=756\", \"386\u00d76=2316\"],\n  [\"942\u00d74=3768\", \"655\u00d79=5895\"],\n  [\"955\u00d76=5730\", \"732\u00d79=6588\"],\n  [\"350\u00d79=3150\", \"655\u00d74=2620\"],\n  [\"913\u00d75=4565\", \"401\u00d73=1203\"],\n  [\"939\u00d74=3756\", \"709\u00d78=5672\"],\n  [\"224\u00d78=1792\", \"353\u00d78=2824\"],\n  [\"998\u00d75=4990\", \"638\u00d72=1276\"],\n  [\"366\u00d77=2562\", \"973\u00d79=8757\"],\n  [\"304\u00d76=1824\", \"632\u00d76=3792\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the three-digit x one-digit multiplication problems/answers\n# in the worksheet table with a newly generated set of problems.\n# Each old equation string is unique in the document, so Find/Replace\n# per pair is safe and precise.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old=\"933\u00d74=3732\"; new=\"425\u00d77=2975\"},\n    @{old=\"738\u00d79=6642\"; new=\"468\u00d79=4212\"},\n    @{old=\"161\u00d74=644\";  new=\"766\u00d77=5362\"},\n    @{old=\"518\u00d76=3108\"; new=\"173\u00d78=1384\"},\n    @{old=\"567\u00d79=5103\"; new=\"101\u00d76=606\"},\n    @{old=\"841\u00d79=7569\"; new=\"173\u00d74=692\"},\n    @{old=\"283\u00d77=1981\"; new=\"562\u00d73=1686\"},\n    @{old=\"359\u00d75=1795\"; new=\"421\u00d75=2105\"},\n    @{old=\"256\u00d74=1024\"; new=\"605\u00d75=3025\"},\n    @{old=\"171\u00d79=1539\"; new=\"291\u00d75=1455\"},\n    @{old=\"997\u00d75=4985\"; new=\"252\u00d79=2268\"},\n    @{old=\"600\u00d72=1200\"; new=\"947\u00d76=5682\"},\n    @{old=\"906\u00d78=7248\"; new=\"897\u00d75=4485\"},\n    @{old=\"281\u00d76=1686\"; new=\"565\u00d79=5085\"},\n    @{old=\"685\u00d79=6165\"; new=\"694\u00d79=6246\"},\n    @{old=\"378\u00d72=756\";  new=\"386\u00d76=2316\"},\n    @{old=\"942\u00d74=3768\"; new=\"655\u00d79=5895\"},\n    @{old=\"955\u00d76=5730\"; new=\"732\u00d79=6588\"},\n    @{old=\"350\u00d79=3150\"; new=\"655\u00d74=2620\"},\n    @{old=\"913\u00d75=4565\"; new=\"401\u00d73=1203\"},\n    @{old=\"939\u00d74=3756\"; new=\"709\u00d78=5672\"},\n    @{old=\"224\u00d78=1792\"; new=\"353\u00d78=2824\"},\n    @{old=\"998\u00d75=4990\"; new=\"638\u00d72=1276\"},\n    @{old=\"366\u00d77=2562\"; new=\"973\u00d79=8757\"},\n    @{old=\"304\u00d76=1824\"; new=\"632\u00d76=3792\"}\n)\n\nforeach ($p in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $p.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $p.new\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n"}
